# RDM-7071: Updated S-578: Added new scenario and data files.
# Adds new "citizen" role authorisation rows for BEFTA_CASETYPE_2_1 across
# the Authorisation* sheets (CaseType, CaseState, CaseEvent, CaseField).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AuthorisationCaseType: new row 9 - citizen gets CRU on BEFTA_CASETYPE_2_1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AuthorisationCaseType")
$ws.Range("A9").Value2 = 42736
$ws.Range("C9").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D9").Value2 = "citizen"
$ws.Range("E9").Value2 = "CRU"

# ---------------------------------------------------------------------------
# AuthorisationCaseState: new rows 19-21 - citizen access to TODO /
# IN_PROGRESS / DONE states of BEFTA_CASETYPE_2_1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AuthorisationCaseState")

$ws.Range("A19").Value2 = 42736
$ws.Range("C19").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D19").Value2 = "TODO"
$ws.Range("E19").Value2 = "citizen"
$ws.Range("F19").Value2 = "CRU"

$ws.Range("A20").Value2 = 42736
$ws.Range("C20").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D20").Value2 = "IN_PROGRESS"
$ws.Range("E20").Value2 = "citizen"
$ws.Range("F20").Value2 = "CRU"

$ws.Range("A21").Value2 = 42736
$ws.Range("C21").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D21").Value2 = "DONE"
$ws.Range("E21").Value2 = "citizen"
$ws.Range("F21").Value2 = "CRU"

# ---------------------------------------------------------------------------
# AuthorisationCaseEvent: fill previously-blank rows 34-39 - citizen access
# to CREATE / START_PROGRESS / STOP_PROGRESS / COMPLETE / UPDATE / REVIEW
# events of BEFTA_CASETYPE_2_1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AuthorisationCaseEvent")

$ws.Range("A34").Value2 = 42736
$ws.Range("C34").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D34").Value2 = "CREATE"
$ws.Range("E34").Value2 = "citizen"
$ws.Range("F34").Value2 = "CRU"

$ws.Range("A35").Value2 = 42736
$ws.Range("C35").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D35").Value2 = "START_PROGRESS"
$ws.Range("E35").Value2 = "citizen"
$ws.Range("F35").Value2 = "CRU"

$ws.Range("A36").Value2 = 42736
$ws.Range("C36").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D36").Value2 = "STOP_PROGRESS"
$ws.Range("E36").Value2 = "citizen"
$ws.Range("F36").Value2 = "CRU"

$ws.Range("A37").Value2 = 42736
$ws.Range("C37").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D37").Value2 = "COMPLETE"
$ws.Range("E37").Value2 = "citizen"
$ws.Range("F37").Value2 = "CRU"

$ws.Range("A38").Value2 = 42736
$ws.Range("C38").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D38").Value2 = "UPDATE"
$ws.Range("E38").Value2 = "citizen"
$ws.Range("F38").Value2 = "CRU"

$ws.Range("A39").Value2 = 42736
$ws.Range("C39").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D39").Value2 = "REVIEW"
$ws.Range("E39").Value2 = "citizen"
$ws.Range("F39").Value2 = "CRU"

# ---------------------------------------------------------------------------
# AuthorisationCaseField: new rows 105-127 - citizen access to every field
# of BEFTA_CASETYPE_2_1 (document fields + base field types)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AuthorisationCaseField")

$ws.Range("A105").Value2 = 42736
$ws.Range("C105").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D105").Value2 = "DocumentField1"
$ws.Range("E105").Value2 = "citizen"
$ws.Range("F105").Value2 = "CRU"

$ws.Range("A106").Value2 = 42736
$ws.Range("C106").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D106").Value2 = "DocumentField2"
$ws.Range("E106").Value2 = "citizen"
$ws.Range("F106").Value2 = "CRU"

$ws.Range("A107").Value2 = 42736
$ws.Range("C107").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D107").Value2 = "DocumentField3"
$ws.Range("E107").Value2 = "citizen"
$ws.Range("F107").Value2 = "CRU"

$ws.Range("A108").Value2 = 42736
$ws.Range("C108").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D108").Value2 = "DocumentField4"
$ws.Range("E108").Value2 = "citizen"
$ws.Range("F108").Value2 = "CR"

$ws.Range("A109").Value2 = 42736
$ws.Range("C109").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D109").Value2 = "DocumentField5"
$ws.Range("E109").Value2 = "citizen"
$ws.Range("F109").Value2 = "CRU"

$ws.Range("A110").Value2 = 42736
$ws.Range("C110").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D110").Value2 = "DocumentField6"
$ws.Range("E110").Value2 = "citizen"
$ws.Range("F110").Value2 = "CR"

$ws.Range("A111").Value2 = 42736
$ws.Range("C111").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D111").Value2 = "DocumentField7"
$ws.Range("E111").Value2 = "citizen"
$ws.Range("F111").Value2 = "CR"

$ws.Range("A112").Value2 = 42736
$ws.Range("C112").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D112").Value2 = "TextField"
$ws.Range("E112").Value2 = "citizen"
$ws.Range("F112").Value2 = "CRU"

$ws.Range("A113").Value2 = 42736
$ws.Range("C113").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D113").Value2 = "NumberField"
$ws.Range("E113").Value2 = "citizen"
$ws.Range("F113").Value2 = "CRU"

$ws.Range("A114").Value2 = 42736
$ws.Range("C114").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D114").Value2 = "YesOrNoField"
$ws.Range("E114").Value2 = "citizen"
$ws.Range("F114").Value2 = "CRU"

$ws.Range("A115").Value2 = 42736
$ws.Range("C115").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D115").Value2 = "PhoneUKField"
$ws.Range("E115").Value2 = "citizen"
$ws.Range("F115").Value2 = "CRU"

$ws.Range("A116").Value2 = 42736
$ws.Range("C116").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D116").Value2 = "EmailField"
$ws.Range("E116").Value2 = "citizen"
$ws.Range("F116").Value2 = "CRU"

$ws.Range("A117").Value2 = 42736
$ws.Range("C117").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D117").Value2 = "MoneyGBPField"
$ws.Range("E117").Value2 = "citizen"
$ws.Range("F117").Value2 = "CRU"

$ws.Range("A118").Value2 = 42736
$ws.Range("C118").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D118").Value2 = "DateField"
$ws.Range("E118").Value2 = "citizen"
$ws.Range("F118").Value2 = "CRU"

$ws.Range("A119").Value2 = 42736
$ws.Range("C119").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D119").Value2 = "DateTimeField"
$ws.Range("E119").Value2 = "citizen"
$ws.Range("F119").Value2 = "CRU"

$ws.Range("A120").Value2 = 42736
$ws.Range("C120").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D120").Value2 = "TextAreaField"
$ws.Range("E120").Value2 = "citizen"
$ws.Range("F120").Value2 = "CRU"

$ws.Range("A121").Value2 = 42736
$ws.Range("C121").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D121").Value2 = "FixedListField"
$ws.Range("E121").Value2 = "citizen"
$ws.Range("F121").Value2 = "CRU"

$ws.Range("A122").Value2 = 42736
$ws.Range("C122").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D122").Value2 = "MultiSelectListField"
$ws.Range("E122").Value2 = "citizen"
$ws.Range("F122").Value2 = "CRU"

$ws.Range("A123").Value2 = 42736
$ws.Range("C123").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D123").Value2 = "AddressUKField"
$ws.Range("E123").Value2 = "citizen"
$ws.Range("F123").Value2 = "CRU"

$ws.Range("A124").Value2 = 42736
$ws.Range("C124").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D124").Value2 = "CollectionField"
$ws.Range("E124").Value2 = "citizen"
$ws.Range("F124").Value2 = "CRU"

$ws.Range("A125").Value2 = 42736
$ws.Range("C125").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D125").Value2 = "ComplexField"
$ws.Range("E125").Value2 = "citizen"
$ws.Range("F125").Value2 = "CRU"

$ws.Range("A126").Value2 = 42736
$ws.Range("C126").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D126").Value2 = "FixedRadioListField"
$ws.Range("E126").Value2 = "citizen"
$ws.Range("F126").Value2 = "CRU"

$ws.Range("A127").Value2 = 42736
$ws.Range("C127").Value2 = "BEFTA_CASETYPE_2_1"
$ws.Range("D127").Value2 = "HistoryComponentField"
$ws.Range("E127").Value2 = "citizen"
$ws.Range("F127").Value2 = "CRU"

# ---------------------------------------------------------------------------
# Leave the final active selection on the AuthorisationComplexType sheet,
# matching the tab that was active when the workbook was last saved.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("AuthorisationComplexType")
$wsFinal.Activate()
$wsFinal.Range("B18").Select()
